$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 830.5  # H34 (was 1288.8)
$ws.Cells.Item(34, 9).Value = 495.75  # I34 (was 722)
$ws.Cells.Item(34, 10).Value = 1500  # J34 (was 1666.6666)
$ws.Cells.Item(34, 11).Value = 495.75  # K34 (was 722)
$ws.Cells.Item(34, 12).Value = 1500  # L34 (was 1666.6666)
$ws.Cells.Item(34, 13).Value = -292.75  # M34 (was -519)
$ws.Cells.Item(34, 14).Value = -1906  # N34 (was -2072.6666)
$ws.Cells.Item(36, 8).Value = 830.5  # H36 (was 1288.8)
$ws.Cells.Item(36, 9).Value = 495.75  # I36 (was 722)
$ws.Cells.Item(36, 10).Value = 1500  # J36 (was 1666.6666)
$ws.Cells.Item(36, 11).Value = 495.75  # K36 (was 722)
$ws.Cells.Item(36, 12).Value = 1500  # L36 (was 1666.6666)
$ws.Cells.Item(36, 13).Value = 219.25  # M36 (was -7)
$ws.Cells.Item(36, 14).Value = -2930  # N36 (was -3096.6666)
$ws.Cells.Item(58, 8).Value = 2804.2  # H58 (was 2919.6667)
$ws.Cells.Item(58, 10).Value = 3492.75  # J58 (was 3493.6)
$ws.Cells.Item(58, 12).Value = 10478.25  # L58 (was 10480.8)
$ws.Cells.Item(58, 14).Value = -10778.25  # N58 (was -10780.8)
$ws.Cells.Item(96, 8).Value = 234  # H96 (was 208.75)
$ws.Cells.Item(96, 9).Value = 262.25  # I96 (was 218.28572)
$ws.Cells.Item(96, 10).Value = 121  # J96 (was 142)
$ws.Cells.Item(96, 11).Value = 786.75  # K96 (was 654.85716)
$ws.Cells.Item(96, 12).Value = 363  # L96 (was 426)
$ws.Cells.Item(96, 13).Value = 586.25  # M96 (was 718.14284)
$ws.Cells.Item(96, 14).Value = -3109  # N96 (was -3172)
$ws.Cells.Item(101, 8).Value = 163  # H101 (was 168)
$ws.Cells.Item(101, 9).Value = 163  # I101 (was 168)
$ws.Cells.Item(101, 11).Value = 489  # K101 (was 504)
$ws.Cells.Item(101, 13).Value = 1133  # M101 (was 1118)
$ws.Cells.Item(113, 8).Value = 8335.333000000001  # H113 (was 6917.5)
$ws.Cells.Item(113, 9).Value = 8000  # I113 (was 6666.6665)
$ws.Cells.Item(113, 10).Value = 8503  # J113 (was 7168.3335)
$ws.Cells.Item(113, 11).Value = 8000  # K113 (was 6666.6665)
$ws.Cells.Item(113, 12).Value = 8503  # L113 (was 7168.3335)
$ws.Cells.Item(113, 13).Value = -4746  # M113 (was -3412.6665)
$ws.Cells.Item(113, 14).Value = -15011  # N113 (was -13676.3335)
$ws.Cells.Item(137, 8).Value = 7132.2856  # H137 (was 6060.857)
$ws.Cells.Item(137, 9).Value = 4700  # I137 (was 2825)
$ws.Cells.Item(137, 11).Value = 14100  # K137 (was 8475)
$ws.Cells.Item(137, 13).Value = -11550  # M137 (was -5925)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9352.333000000001  # H32 (was 1497)
$ws.Cells.Item(32, 9).Value = 3278.5  # I32 (was 1121.25)
$ws.Cells.Item(32, 10).Value = 21500  # J32 (was 3000)
$ws.Cells.Item(32, 11).Value = 3278.5  # K32 (was 1121.25)
$ws.Cells.Item(32, 12).Value = 21500  # L32 (was 3000)
$ws.Cells.Item(32, 13).Value = -2991.5  # M32 (was -834.25)
$ws.Cells.Item(32, 14).Value = -22074  # N32 (was -3574)
$ws.Cells.Item(43, 8).Value = 49990  # H43 (was 49683.5)
$ws.Cells.Item(43, 10).Value = 49990  # J43 (was 49683.5)
$ws.Cells.Item(43, 12).Value = 49990  # L43 (was 49683.5)
$ws.Cells.Item(43, 14).Value = -50616  # N43 (was -50309.5)
$ws.Cells.Item(45, 8).Value = 7879.2354  # H45 (was 8663.134)
$ws.Cells.Item(45, 9).Value = 4474.9165  # I45 (was 5299.8887)
$ws.Cells.Item(45, 10).Value = 16049.6  # J45 (was 13708)
$ws.Cells.Item(45, 11).Value = 4474.9165  # K45 (was 5299.8887)
$ws.Cells.Item(45, 12).Value = 16049.6  # L45 (was 13708)
$ws.Cells.Item(45, 13).Value = -4097.9165  # M45 (was -4922.8887)
$ws.Cells.Item(45, 14).Value = -16803.6  # N45 (was -14462)
$ws.Cells.Item(97, 8).Value = 1397.8  # H97 (was 1502.25)
$ws.Cells.Item(97, 9).Value = 1497.25  # I97 (was 1669.6666)
$ws.Cells.Item(97, 11).Value = 1497.25  # K97 (was 1669.6666)
$ws.Cells.Item(97, 13).Value = -1001.25  # M97 (was -1173.6666)
$ws.Cells.Item(102, 8).Value = 4737  # H102 (was 3634.25)
$ws.Cells.Item(102, 9).Value = 4649.3335  # I102 (was 3439.1428)
$ws.Cells.Item(102, 11).Value = 4649.3335  # K102 (was 3439.1428)
$ws.Cells.Item(102, 13).Value = -3027.3335  # M102 (was -1817.1428)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1953.6666  # H107 (was 2749.6667)
$ws.Cells.Item(107, 9).Value = 1572.875  # I107 (was 2299.6)
$ws.Cells.Item(107, 11).Value = 1572.875  # K107 (was 2299.6)
$ws.Cells.Item(107, 13).Value = 347.125  # M107 (was -379.5999999999999)
$ws.Cells.Item(135, 8).Value = 0  # H135 (was 60000)
$ws.Cells.Item(135, 10).Value = 0  # J135 (was 60000)
$ws.Cells.Item(135, 12).Value = 0  # L135 (was 60000)
$ws.Cells.Item(135, 14).Value = ""  # N135 (was -70140)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 461.8  # H22 (was 1079.4)
$ws.Cells.Item(22, 9).Value = 413.36365  # I22 (was 1598.5)
$ws.Cells.Item(22, 10).Value = 595  # J22 (was 733.3333)
$ws.Cells.Item(22, 11).Value = 413.36365  # K22 (was 1598.5)
$ws.Cells.Item(22, 12).Value = 595  # L22 (was 733.3333)
$ws.Cells.Item(22, 13).Value = -63.36365000000001  # M22 (was -1248.5)
$ws.Cells.Item(22, 14).Value = -1295  # N22 (was -1433.3333)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 201.25  # H38 (was 137.3)
$ws.Cells.Item(38, 10).Value = 306.5  # J38 (was 199.25)
$ws.Cells.Item(38, 12).Value = 919.5  # L38 (was 597.75)
$ws.Cells.Item(38, 14).Value = -1613.5  # N38 (was -1291.75)
$ws.Cells.Item(80, 8).Value = 5000  # H80 (was 7900)
$ws.Cells.Item(80, 9).Value = 5000  # I80 (was 0)
$ws.Cells.Item(80, 10).Value = 0  # J80 (was 7900)
$ws.Cells.Item(80, 11).Value = 15000  # K80 (was 0)
$ws.Cells.Item(80, 12).Value = 0  # L80 (was 23700)
$ws.Cells.Item(80, 13).Value = -14064  # M80 (was None)
$ws.Cells.Item(80, 14).Value = ""  # N80 (was -25572)
$ws.Cells.Item(83, 8).Value = 5000  # H83 (was 7900)
$ws.Cells.Item(83, 9).Value = 5000  # I83 (was 0)
$ws.Cells.Item(83, 10).Value = 0  # J83 (was 7900)
$ws.Cells.Item(83, 11).Value = 45000  # K83 (was 0)
$ws.Cells.Item(83, 12).Value = 0  # L83 (was 71100)
$ws.Cells.Item(83, 13).Value = -40320  # M83 (was None)
$ws.Cells.Item(83, 14).Value = ""  # N83 (was -80460)
$ws.Cells.Item(122, 8).Value = 475.22223  # H122 (was 387.25)
$ws.Cells.Item(122, 9).Value = 300  # I122 (was 208.33333)
$ws.Cells.Item(122, 10).Value = 562.8333  # J122 (was 566.1667)
$ws.Cells.Item(122, 11).Value = 2700  # K122 (was 1874.99997)
$ws.Cells.Item(122, 12).Value = 5065.4997  # L122 (was 5095.5003)
$ws.Cells.Item(122, 13).Value = -250  # M122 (was 575.0000300000002)
$ws.Cells.Item(122, 14).Value = -9965.4997  # N122 (was -9995.5003)
$ws.Cells.Item(129, 8).Value = 862.6  # H129 (was 1026.75)
$ws.Cells.Item(129, 9).Value = 578.25  # I129 (was 702.3333)
$ws.Cells.Item(129, 11).Value = 1734.75  # K129 (was 2106.9999)
$ws.Cells.Item(129, 13).Value = 3265.25  # M129 (was 2893.0001)
$ws.Cells.Item(131, 8).Value = 1809.7  # H131 (was 1899.6666)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 1000000  # H18 (was 410)
$ws.Cells.Item(18, 9).Value = 1000000  # I18 (was 410)
$ws.Cells.Item(18, 11).Value = 1000000  # K18 (was 410)
$ws.Cells.Item(18, 13).Value = -999707  # M18 (was -117)
$ws.Cells.Item(70, 8).Value = 5333.3335  # H70 (was 5000)
$ws.Cells.Item(70, 10).Value = 6000  # J70 (was 0)
$ws.Cells.Item(70, 12).Value = 6000  # L70 (was 0)
$ws.Cells.Item(70, 14).Value = -6540  # N70 (was None)
$ws.Cells.Item(73, 8).Value = 5333.3335  # H73 (was 5000)
$ws.Cells.Item(73, 10).Value = 6000  # J73 (was 0)
$ws.Cells.Item(73, 12).Value = 6000  # L73 (was 0)
$ws.Cells.Item(73, 14).Value = -7872  # N73 (was None)
$ws.Cells.Item(103, 8).Value = 40000  # H103 (was 0)
$ws.Cells.Item(103, 10).Value = 40000  # J103 (was 0)
$ws.Cells.Item(103, 12).Value = 40000  # L103 (was 0)
$ws.Cells.Item(103, 14).Value = -42344  # N103 (was None)
$ws.Cells.Item(122, 8).Value = 3787.75  # H122 (was 5001.1665)
$ws.Cells.Item(122, 9).Value = 3787.75  # I122 (was 5001.1665)
$ws.Cells.Item(122, 11).Value = 11363.25  # K122 (was 15003.4995)
$ws.Cells.Item(122, 13).Value = -8913.25  # M122 (was -12553.4995)
$ws.Cells.Item(132, 8).Value = 6146.8335  # H132 (was 6600.222)
$ws.Cells.Item(132, 9).Value = 2977.05  # I132 (was 3520.6667)
$ws.Cells.Item(132, 10).Value = 21995.75  # J132 (was 21998)
$ws.Cells.Item(132, 11).Value = 8931.150000000001  # K132 (was 10562.0001)
$ws.Cells.Item(132, 12).Value = 65987.25  # L132 (was 65994)
$ws.Cells.Item(132, 13).Value = -6401.150000000001  # M132 (was -8032.000100000001)
$ws.Cells.Item(132, 14).Value = -71047.25  # N132 (was -71054)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(56, 8).Value = 100051  # H56 (was 55525)
$ws.Cells.Item(56, 9).Value = 100051  # I56 (was 55525)
$ws.Cells.Item(56, 11).Value = 100051  # K56 (was 55525)
$ws.Cells.Item(56, 13).Value = -99360  # M56 (was -54834)
$ws.Cells.Item(93, 8).Value = 4285.4287  # H93 (was 4499.8335)
$ws.Cells.Item(93, 10).Value = 2999  # J93 (was 0)
$ws.Cells.Item(93, 12).Value = 2999  # L93 (was 0)
$ws.Cells.Item(93, 14).Value = -5495  # N93 (was None)
$ws.Cells.Item(122, 8).Value = 3969.125  # H122 (was 4025.75)
$ws.Cells.Item(122, 9).Value = 3851.3333  # I122 (was 3904)
$ws.Cells.Item(122, 10).Value = 4039.8  # J122 (was 4066.3333)
$ws.Cells.Item(122, 11).Value = 11553.9999  # K122 (was 11712)
$ws.Cells.Item(122, 12).Value = 12119.4  # L122 (was 12198.9999)
$ws.Cells.Item(122, 13).Value = -9103.999899999999  # M122 (was -9262)
$ws.Cells.Item(122, 14).Value = -17019.4  # N122 (was -17098.9999)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17, 8).Value = 999  # H17 (was 0)
$ws.Cells.Item(17, 10).Value = 999  # J17 (was 0)
$ws.Cells.Item(17, 12).Value = 999  # L17 (was 0)
$ws.Cells.Item(17, 14).Value = -1343  # N17 (was None)
$ws.Cells.Item(113, 8).Value = 633.8889  # H113 (was 780.8333)
$ws.Cells.Item(113, 9).Value = 620.5714  # I113 (was 780.8333)
$ws.Cells.Item(113, 10).Value = 680.5  # J113 (was 0)
$ws.Cells.Item(113, 11).Value = 1861.7142  # K113 (was 2342.4999)
$ws.Cells.Item(113, 12).Value = 2041.5  # L113 (was 0)
$ws.Cells.Item(113, 13).Value = 308.2857999999999  # M113 (was -172.4998999999998)
$ws.Cells.Item(113, 14).Value = -6381.5  # N113 (was None)
$ws.Cells.Item(130, 8).Value = 29991  # H130 (was 29995)
$ws.Cells.Item(130, 10).Value = 29991  # J130 (was 29995)
$ws.Cells.Item(130, 12).Value = 29991  # L130 (was 29995)
$ws.Cells.Item(130, 14).Value = -40031  # N130 (was -40035)
